$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume table with the latest snapshot.
# Price cells (column D) sometimes look numeric (e.g. "0.999", "69.00",
# "172.30") but must stay literal text so formatting like trailing
# zeros and thousand-separator dots ("65.687.81") survives -- briefly
# force a text NumberFormat while assigning, then restore the default
# "Normal" style so no stray per-cell formatting is left behind.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "65.687.81"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.58%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.167.36"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -4.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "572.28"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "172.30"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.11%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "3.166.85"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("E11").Value = "  -3.52%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.392"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.42%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.718.20"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -4.89%  "
$ws.Range("E14").Value = "  +0.98%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.27"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.53%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.677.01"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  -2.30%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.171.00"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.34%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.91"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.15%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "361.11"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  -1.63%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "69.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("E25").Value = "  -4.40%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.307.90"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.06%  "
$ws.Range("E28").Value = "  +3.83%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -0.11%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.41"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "22.07"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.29%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.64"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("E36").Value = "  -0.65%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "160.78"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("E40").Value = "  +2.95%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.34"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("E42").Value = "  -2.46%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.651.23"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("E44").Value = "  -0.28%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.19"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.43%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "39.74"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0657"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "330.12"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.79%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "23.87"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0275"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -0.62%  "
